$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Rename "Rule 4" -> "Project Units"
# ------------------------------------------------------------------
$wsUnits = $wb.Worksheets.Item("Rule 4")
$wsUnits.Name = "Project Units"

# ------------------------------------------------------------------
# 2) "Not Placed" sheet: merge header row, center it, move selection
# ------------------------------------------------------------------
$wsNotPlaced = $wb.Worksheets.Item("Not Placed")
$wsNotPlaced.Range("A1:F1").HorizontalAlignment = -4108
$wsNotPlaced.Range("A1:F1").Merge()
$wsNotPlaced.Range("F4").Select()

# ------------------------------------------------------------------
# 3) Build out the "Project Units" sheet contents
# ------------------------------------------------------------------
$wsUnits.Range("A1").Value = "Rule check results (Project Units)"
$wsUnits.Range("A1:O1").HorizontalAlignment = -4108
$wsUnits.Range("A1:O1").Merge()

$headers = @("Length (Units)","Length (Rounding)","Area (Units)","Area (Rounding)","Volume (Units)","Volume (Rounding)","Angle (Units)","Angle (Rounding)","Slope (Units)","Slope (Rounding)","Currency (Units)","Currency (Rounding)","Mass Density (Units)","Mass Density (Rounding)","File")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsUnits.Cells.Item(3, $i + 1).Value = $headers[$i]
}

$values = @("Milimeters","0 decimal places","Square meters","2 decimal places","Cubic meters","2 decimal meters","Decimal degrees","2 decimal places","Percentage","2 decimal places","Currency","2 decimal places","Kilograms per cubic meter","2 decimal places","Z:\02\Proyectos\Habana….")
for ($i = 0; $i -lt $values.Length; $i++) {
    $wsUnits.Cells.Item(4, $i + 1).Value = $values[$i]
}

$tableRange = $wsUnits.Range("A3:O4")
$lo = $wsUnits.ListObjects.Add(1, $tableRange, 0, 1)
$lo.Name = "Tabla2"
$lo.TableStyle = "TableStyleLight1"

$wsUnits.Range("A3:O3").Select()
$wsUnits.Activate()
$excel.ActiveWindow.ScrollColumn = 3

# ------------------------------------------------------------------
# 4) Set HOME sheet active tab selection back to "Project Units"
# ------------------------------------------------------------------
$wsUnits.Activate()
